$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The sheet currently ends at row 84 (the most recent Región del Maule
# record). Two newer weekly records need to be inserted above it, pushing
# the existing row 84 down to row 86.
$ws.Rows("84:85").Insert()

# New row 84: "Especial" quality record for 2022-02-03 (serial 44595),
# origin "Provincia de Linares".
$ws.Range("A84").Value = 9
$ws.Range("B84").Value = "Vega Central Mapocho de Santiago"
$ws.Range("C84").Value = "Metropolitana"
$ws.Range("D84").Value = 44595
$ws.Range("E84").Value = 13
$ws.Range("F84").Value = "Fruta"
$ws.Range("G84").Value = 100101
$ws.Range("H84").Value = "Berries"
$ws.Range("I84").Value = 100101004
$ws.Range("J84").Value = "Frambuesa"
$ws.Range("K84").Value = "Sin especificar"
$ws.Range("L84").Value = "Especial"
$ws.Range("M84").Value = 400
$ws.Range("N84").Value = 8000
$ws.Range("O84").Value = 8000
$ws.Range("P84").Value = 8000
$ws.Range("Q84").Value = "$/bandeja 2 kilos"
$ws.Range("R84").Value = "Provincia de Linares"
$ws.Range("S84").Value = 4000
$ws.Range("T84").Value = 2

# New row 85: "Primera" quality record, same date, same origin.
$ws.Range("A85").Value = 9
$ws.Range("B85").Value = "Vega Central Mapocho de Santiago"
$ws.Range("C85").Value = "Metropolitana"
$ws.Range("D85").Value = 44595
$ws.Range("E85").Value = 13
$ws.Range("F85").Value = "Fruta"
$ws.Range("G85").Value = 100101
$ws.Range("H85").Value = "Berries"
$ws.Range("I85").Value = 100101004
$ws.Range("J85").Value = "Frambuesa"
$ws.Range("K85").Value = "Sin especificar"
$ws.Range("L85").Value = "Primera"
$ws.Range("M85").Value = 450
$ws.Range("N85").Value = 7000
$ws.Range("O85").Value = 7000
$ws.Range("P85").Value = 7000
$ws.Range("Q85").Value = "$/bandeja 2 kilos"
$ws.Range("R85").Value = "Provincia de Linares"
$ws.Range("S85").Value = 3500
$ws.Range("T85").Value = 2

# Row 86 already holds the original row-84 data (shifted down by the
# insert above), so nothing else to do there.
